$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Text
    if ($val -eq "MOLLY MCNINCH") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENT") {
        $cell.Value = "S"
    }
}
